$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("drivers_data")

# --- Row 37: LIO CAR REPAIRS LTD ---
$ws.Range("C2").Copy()
$ws.Range("A37").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A37").HorizontalAlignment = -4108 # xlCenter
$ws.Range("A37").VerticalAlignment = -4108 # xlCenter

$ws.Range("B25").Copy()
$ws.Range("B37").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("C2").Copy()
$ws.Range("C37").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("D25:H25").Copy()
$ws.Range("D37:H37").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A37").Value = "LIO CAR REPAIRS LTD"
$ws.Range("B37").Value = "07491 110752"
$ws.Range("C37").Value = "101 Roosevelt Way, Dagenham, United Kingdom, RM10 8DA"
$ws.Range("D37").Value = "Yes"
$ws.Range("E37").Value = "No"
$ws.Range("F37").Value = "No"
$ws.Range("G37").Value = "No"
$ws.Range("H37").Value = "No"

# --- Row 38: ASL CAR MOTORCYCLE BREAKDOWN RECOVERY ---
$ws.Range("C2").Copy()
$ws.Range("A38").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A38").HorizontalAlignment = -4108 # xlCenter
$ws.Range("A38").VerticalAlignment = -4108 # xlCenter

$ws.Range("B25").Copy()
$ws.Range("B38").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("C2").Copy()
$ws.Range("C38").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("D25:H25").Copy()
$ws.Range("D38:H38").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A38").Value = "ASL CAR MOTORCYCLE BREAKDOWN RECOVERY"
$ws.Range("B38").Value = "07950 355777"
$ws.Range("C38").Value = "56 tavistock road 165 epsom road London SM5 1QR"
$ws.Range("D38").Value = "Yes"
$ws.Range("E38").Value = "No"
$ws.Range("F38").Value = "No"
$ws.Range("G38").Value = "No"
$ws.Range("H38").Value = "No"

$ws.Rows.Item(37).RowHeight = 15.75
$ws.Rows.Item(38).RowHeight = 15.75

$excel.CutCopyMode = 0

# Update the view to match target: topLeftCell A16, selection F42
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("F42").Select() | Out-Null
